# Apply "2020-09-02 data" update to the Fonds de solidarite volet-1 dataset.
# Updates nombre_aides (column C) and montant_total (column D) for the rows
# that received new cumulative figures in this data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

    $ws.Cells.Item(2, 3).Value = 38885
    $ws.Cells.Item(2, 4).Value = 56221263
    $ws.Cells.Item(3, 3).Value = 93209
    $ws.Cells.Item(3, 4).Value = 136620068
    $ws.Cells.Item(4, 3).Value = 31830
    $ws.Cells.Item(4, 4).Value = 47133967
    $ws.Cells.Item(5, 3).Value = 8946
    $ws.Cells.Item(5, 4).Value = 13296397
    $ws.Cells.Item(6, 3).Value = 2082
    $ws.Cells.Item(6, 4).Value = 3094971
    $ws.Cells.Item(12, 3).Value = 42285
    $ws.Cells.Item(12, 4).Value = 57322564
    $ws.Cells.Item(13, 3).Value = 9917
    $ws.Cells.Item(13, 4).Value = 14340528
    $ws.Cells.Item(14, 3).Value = 26498
    $ws.Cells.Item(14, 4).Value = 38842826
    $ws.Cells.Item(15, 3).Value = 8472
    $ws.Cells.Item(15, 4).Value = 12573478
    $ws.Cells.Item(16, 3).Value = 2211
    $ws.Cells.Item(16, 4).Value = 3285039
    $ws.Cells.Item(17, 3).Value = 431
    $ws.Cells.Item(17, 4).Value = 635623
    $ws.Cells.Item(18, 3).Value = 37
    $ws.Cells.Item(18, 4).Value = 55500
    $ws.Cells.Item(20, 3).Value = 10426
    $ws.Cells.Item(20, 4).Value = 13791047
    $ws.Cells.Item(21, 3).Value = 13716
    $ws.Cells.Item(21, 4).Value = 19794066
    $ws.Cells.Item(22, 3).Value = 32280
    $ws.Cells.Item(22, 4).Value = 47362423
    $ws.Cells.Item(23, 3).Value = 10432
    $ws.Cells.Item(23, 4).Value = 15505548
    $ws.Cells.Item(24, 3).Value = 2702
    $ws.Cells.Item(24, 4).Value = 4017771
    $ws.Cells.Item(25, 3).Value = 540
    $ws.Cells.Item(25, 4).Value = 804092
    $ws.Cells.Item(27, 3).Value = 11936
    $ws.Cells.Item(27, 4).Value = 15932775
    $ws.Cells.Item(28, 3).Value = 7885
    $ws.Cells.Item(28, 4).Value = 11409877
    $ws.Cells.Item(29, 3).Value = 23029
    $ws.Cells.Item(29, 4).Value = 33805558
    $ws.Cells.Item(30, 3).Value = 7963
    $ws.Cells.Item(30, 4).Value = 11843392
    $ws.Cells.Item(31, 3).Value = 2014
    $ws.Cells.Item(31, 4).Value = 3005251
    $ws.Cells.Item(34, 3).Value = 8492
    $ws.Cells.Item(34, 4).Value = 11215451
    $ws.Cells.Item(35, 3).Value = 3361
    $ws.Cells.Item(35, 4).Value = 4854691
    $ws.Cells.Item(36, 3).Value = 8037
    $ws.Cells.Item(36, 4).Value = 11737926
    $ws.Cells.Item(37, 3).Value = 3236
    $ws.Cells.Item(37, 4).Value = 4796961
    $ws.Cells.Item(38, 3).Value = 841
    $ws.Cells.Item(38, 4).Value = 1252723
    $ws.Cells.Item(40, 3).Value = 5
    $ws.Cells.Item(40, 4).Value = 7500
    $ws.Cells.Item(41, 3).Value = 2535
    $ws.Cells.Item(41, 4).Value = 3424349
    $ws.Cells.Item(42, 3).Value = 17758
    $ws.Cells.Item(42, 4).Value = 25679057
    $ws.Cells.Item(43, 3).Value = 52280
    $ws.Cells.Item(43, 4).Value = 76625271
    $ws.Cells.Item(44, 3).Value = 19344
    $ws.Cells.Item(44, 4).Value = 28726466
    $ws.Cells.Item(45, 3).Value = 5743
    $ws.Cells.Item(45, 4).Value = 8548187
    $ws.Cells.Item(46, 3).Value = 1260
    $ws.Cells.Item(46, 4).Value = 1880545
    $ws.Cells.Item(47, 3).Value = 69
    $ws.Cells.Item(47, 4).Value = 101568
    $ws.Cells.Item(50, 3).Value = 17161
    $ws.Cells.Item(50, 4).Value = 22796694
    $ws.Cells.Item(51, 3).Value = 2150
    $ws.Cells.Item(51, 4).Value = 3120878
    $ws.Cells.Item(52, 3).Value = 7237
    $ws.Cells.Item(52, 4).Value = 10635725
    $ws.Cells.Item(53, 3).Value = 2433
    $ws.Cells.Item(53, 4).Value = 3633572
    $ws.Cells.Item(54, 3).Value = 771
    $ws.Cells.Item(54, 4).Value = 1151805
    $ws.Cells.Item(55, 3).Value = 199
    $ws.Cells.Item(55, 4).Value = 294726
    $ws.Cells.Item(57, 3).Value = 7355
    $ws.Cells.Item(57, 4).Value = 10115160
    $ws.Cells.Item(58, 3).Value = 1183
    $ws.Cells.Item(58, 4).Value = 2017793
    $ws.Cells.Item(59, 3).Value = 2889
    $ws.Cells.Item(59, 4).Value = 4908753
    $ws.Cells.Item(60, 3).Value = 1133
    $ws.Cells.Item(60, 4).Value = 1927026
    $ws.Cells.Item(61, 3).Value = 395
    $ws.Cells.Item(61, 4).Value = 681883
    $ws.Cells.Item(62, 3).Value = 128
    $ws.Cells.Item(62, 4).Value = 223100
    $ws.Cells.Item(64, 3).Value = 1719
    $ws.Cells.Item(64, 4).Value = 2700578
    $ws.Cells.Item(65, 3).Value = 15834
    $ws.Cells.Item(65, 4).Value = 22869285
    $ws.Cells.Item(66, 3).Value = 45781
    $ws.Cells.Item(66, 4).Value = 66977705
    $ws.Cells.Item(67, 3).Value = 16000
    $ws.Cells.Item(67, 4).Value = 23772506
    $ws.Cells.Item(68, 3).Value = 4660
    $ws.Cells.Item(68, 4).Value = 6939788
    $ws.Cells.Item(69, 3).Value = 968
    $ws.Cells.Item(69, 4).Value = 1440168
    $ws.Cells.Item(70, 3).Value = 80
    $ws.Cells.Item(70, 4).Value = 117330
    $ws.Cells.Item(73, 3).Value = 15435
    $ws.Cells.Item(73, 4).Value = 20327382
    $ws.Cells.Item(74, 3).Value = 54719
    $ws.Cells.Item(74, 4).Value = 79623105
    $ws.Cells.Item(75, 3).Value = 152858
    $ws.Cells.Item(75, 4).Value = 225169543
    $ws.Cells.Item(76, 3).Value = 65852
    $ws.Cells.Item(76, 4).Value = 98121147
    $ws.Cells.Item(77, 3).Value = 21120
    $ws.Cells.Item(77, 4).Value = 31557884
    $ws.Cells.Item(78, 3).Value = 5052
    $ws.Cells.Item(78, 4).Value = 7546403
    $ws.Cells.Item(79, 3).Value = 284
    $ws.Cells.Item(79, 4).Value = 421170
    $ws.Cells.Item(80, 3).Value = 24
    $ws.Cells.Item(80, 4).Value = 34905
    $ws.Cells.Item(85, 3).Value = 53710
    $ws.Cells.Item(85, 4).Value = 72960206
    $ws.Cells.Item(86, 3).Value = 4782
    $ws.Cells.Item(86, 4).Value = 6929220
    $ws.Cells.Item(87, 3).Value = 11908
    $ws.Cells.Item(87, 4).Value = 17489562
    $ws.Cells.Item(88, 3).Value = 3969
    $ws.Cells.Item(88, 4).Value = 5913458
    $ws.Cells.Item(89, 3).Value = 1371
    $ws.Cells.Item(89, 4).Value = 2048289
    $ws.Cells.Item(90, 3).Value = 296
    $ws.Cells.Item(90, 4).Value = 441512
    $ws.Cells.Item(93, 3).Value = 5571
    $ws.Cells.Item(93, 4).Value = 7486502
    $ws.Cells.Item(94, 3).Value = 1658
    $ws.Cells.Item(94, 4).Value = 2389199
    $ws.Cells.Item(95, 3).Value = 5346
    $ws.Cells.Item(95, 4).Value = 7874519
    $ws.Cells.Item(96, 3).Value = 1987
    $ws.Cells.Item(96, 4).Value = 2958426
    $ws.Cells.Item(97, 3).Value = 706
    $ws.Cells.Item(97, 4).Value = 1057960
    $ws.Cells.Item(101, 3).Value = 3704
    $ws.Cells.Item(101, 4).Value = 4906713
    $ws.Cells.Item(102, 3).Value = 729
    $ws.Cells.Item(102, 4).Value = 1229775
    $ws.Cells.Item(103, 3).Value = 456
    $ws.Cells.Item(103, 4).Value = 804027
    $ws.Cells.Item(104, 3).Value = 167
    $ws.Cells.Item(104, 4).Value = 291180
    $ws.Cells.Item(105, 3).Value = 54
    $ws.Cells.Item(105, 4).Value = 93000
    $ws.Cells.Item(106, 3).Value = 29
    $ws.Cells.Item(106, 4).Value = 55500
    $ws.Cells.Item(107, 3).Value = 11097
    $ws.Cells.Item(107, 4).Value = 16095029
    $ws.Cells.Item(108, 3).Value = 29800
    $ws.Cells.Item(108, 4).Value = 43762178
    $ws.Cells.Item(109, 3).Value = 9984
    $ws.Cells.Item(109, 4).Value = 14845205
    $ws.Cells.Item(110, 3).Value = 2757
    $ws.Cells.Item(110, 4).Value = 4110580
    $ws.Cells.Item(111, 3).Value = 508
    $ws.Cells.Item(111, 4).Value = 757046
    $ws.Cells.Item(112, 3).Value = 55
    $ws.Cells.Item(112, 4).Value = 82500
    $ws.Cells.Item(114, 3).Value = 10001
    $ws.Cells.Item(114, 4).Value = 13205265
    $ws.Cells.Item(115, 3).Value = 31268
    $ws.Cells.Item(115, 4).Value = 45078066
    $ws.Cells.Item(116, 3).Value = 67569
    $ws.Cells.Item(116, 4).Value = 98861852
    $ws.Cells.Item(117, 3).Value = 21774
    $ws.Cells.Item(117, 4).Value = 32354288
    $ws.Cells.Item(118, 3).Value = 6188
    $ws.Cells.Item(118, 4).Value = 9217841
    $ws.Cells.Item(119, 3).Value = 1161
    $ws.Cells.Item(119, 4).Value = 1735100
    $ws.Cells.Item(120, 3).Value = 86
    $ws.Cells.Item(120, 4).Value = 125395
    $ws.Cells.Item(124, 3).Value = 26361
    $ws.Cells.Item(124, 4).Value = 35174041
    $ws.Cells.Item(125, 3).Value = 37102
    $ws.Cells.Item(125, 4).Value = 53532216
    $ws.Cells.Item(126, 3).Value = 78675
    $ws.Cells.Item(126, 4).Value = 115027508
    $ws.Cells.Item(127, 3).Value = 24326
    $ws.Cells.Item(127, 4).Value = 36104541
    $ws.Cells.Item(128, 3).Value = 6540
    $ws.Cells.Item(128, 4).Value = 9718623
    $ws.Cells.Item(129, 3).Value = 1299
    $ws.Cells.Item(129, 4).Value = 1931811
    $ws.Cells.Item(130, 3).Value = 67
    $ws.Cells.Item(130, 4).Value = 98728
    $ws.Cells.Item(131, 3).Value = 19
    $ws.Cells.Item(131, 4).Value = 28500
    $ws.Cells.Item(133, 3).Value = 32510
    $ws.Cells.Item(133, 4).Value = 43143397
    $ws.Cells.Item(134, 3).Value = 13638
    $ws.Cells.Item(134, 4).Value = 19741797
    $ws.Cells.Item(135, 3).Value = 33020
    $ws.Cells.Item(135, 4).Value = 48491090
    $ws.Cells.Item(136, 3).Value = 11695
    $ws.Cells.Item(136, 4).Value = 17376587
    $ws.Cells.Item(137, 3).Value = 3033
    $ws.Cells.Item(137, 4).Value = 4520741
    $ws.Cells.Item(139, 3).Value = 37
    $ws.Cells.Item(139, 4).Value = 54325
    $ws.Cells.Item(141, 3).Value = 11035
    $ws.Cells.Item(141, 4).Value = 14705717
    $ws.Cells.Item(142, 3).Value = 36200
    $ws.Cells.Item(142, 4).Value = 52278755
    $ws.Cells.Item(143, 3).Value = 83538
    $ws.Cells.Item(143, 4).Value = 122382876
    $ws.Cells.Item(144, 3).Value = 24900
    $ws.Cells.Item(144, 4).Value = 36990558
    $ws.Cells.Item(145, 3).Value = 6542
    $ws.Cells.Item(145, 4).Value = 9761496
    $ws.Cells.Item(146, 3).Value = 1487
    $ws.Cells.Item(146, 4).Value = 2212730
    $ws.Cells.Item(147, 3).Value = 86
    $ws.Cells.Item(147, 4).Value = 128630
    $ws.Cells.Item(149, 3).Value = 29895
    $ws.Cells.Item(149, 4).Value = 40301416
